$wb = $excel.ActiveWorkbook

# --- Window view tweaks (best-effort; matches target bookView numbers) ---
try {
    $win = $wb.Windows.Item(1)
    $win.Left = -110
    $win.Top = -110
    $win.Width = 19420
    $win.Height = 10300
} catch {}

# --- Add the new "Web Shear" worksheet at the very end of the workbook ---
# A throwaway sheet is created first (and deleted afterwards) purely so the
# internal sheetId counter advances past the one that would otherwise be
# reused; this reproduces sheetId="9" on the real sheet instead of "8".
$placeholder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTmp.Name = "Web Shear"
$placeholder.Delete()

# Re-acquire the worksheet by name: the reference obtained before the
# placeholder's deletion can point at a stale sheet-index internally.
$ws = $wb.Worksheets.Item("Web Shear")
$ws.Activate()

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 13.54296875
$ws.Columns.Item(5).ColumnWidth = 22.81640625
$ws.Columns.Item(6).ColumnWidth = 26.26953125
$ws.Columns.Item(7).ColumnWidth = 12.08984375

# --- Header row ---
$ws.Range("A1").Value = "Elastic modulus MPa"
$ws.Range("B1").Value = "yield stress MPa"
$ws.Range("C1").Value = "web area mm"
$ws.Range("D1").Value = "web ratio"
$ws.Range("E1").Value = "rolled web shear limit "
$ws.Range("F1").Value = "web shear strength coefficient"
$ws.Range("G1").Value = "nominal strength"
$ws.Range("H1").Value = "design strength"

# --- Row 2 data ---
$ws.Range("A2").Value = 200000
$ws.Range("B2").Value = 355
$ws.Range("C2").Value = 29344
$ws.Range("D2").Value = 38
$ws.Range("E2").Formula = "=2.24*(A2/B2)^0.5"
$ws.Range("F2").Value = 1
$ws.Range("G2").Formula = "=0.6*B2*C2*F2"
$ws.Range("H2").Formula = "=G2/1.5"

# --- Row 3 data ---
$ws.Range("A3").Value = 200000
$ws.Range("B3").Value = 250
$ws.Range("C3").Value = 887.68
$ws.Range("D3").Value = 21.6
$ws.Range("E3").Formula = "=2.24*(A3/B3)^0.5"
$ws.Range("F3").Value = 1
$ws.Range("G3").Formula = "=0.6*B3*C3*F3"
$ws.Range("H3").Formula = "=G3"

$ws.Range("H4").Select()
